# Refresh the crypto price/volume table (and a handful of shifted
# Coin/Link rows) to match the latest scrape, per commit
# "Updated symbol list on Wed Feb  8 23:48:59 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''328.42'

$ws.Range("D3").Value = '''45.10'
$ws.Range("E3").Value = '''-1.29%'

$ws.Range("D4").Value = '''5.213'
$ws.Range("E4").Value = '''-6.29%'

$ws.Range("D5").Value = '''0.08376'
$ws.Range("E5").Value = '''0.48%'

$ws.Range("D6").Value = '''1.960'
$ws.Range("E6").Value = '''-4.25%'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''4.429'
$ws.Range("E7").Value = '''-0.14%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9733'
$ws.Range("E8").Value = '''-0.78%'

$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '''2.532'
$ws.Range("E9").Value = '''-3.76%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1105'
$ws.Range("E10").Value = '''-2.96%'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1915'
$ws.Range("E11").Value = '''-1.28%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09690'
$ws.Range("E12").Value = '''-4.03%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04608'
$ws.Range("E13").Value = '''-0.51%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1062'
$ws.Range("E14").Value = '''0.03%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001293'
$ws.Range("E15").Value = '''1.88%'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005805'
$ws.Range("E16").Value = '''-3.70%'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.366'
$ws.Range("E17").Value = '''0.03%'

$ws.Range("E18").Value = '''0.21%'

$ws.Range("D19").Value = '''8.367'
$ws.Range("E19").Value = '''-18.35%'

$ws.Range("D20").Value = '''0.1352'

$ws.Range("E21").Value = '''9.18%'

$ws.Range("D22").Value = '''0.04176'
$ws.Range("E22").Value = '''1.50%'

$ws.Range("D23").Value = '''0.001238'
$ws.Range("E23").Value = '''-4.95%'

$ws.Range("D24").Value = '''0.004457'
$ws.Range("E24").Value = '''0.78%'

$ws.Range("D25").Value = '''0.0001300'
$ws.Range("E25").Value = '''1.70%'

$ws.Range("D26").Value = '''0.0002980'
$ws.Range("E26").Value = '''-20.29%'

$ws.Range("D38").Value = '''0.02715'
$ws.Range("E38").Value = '''-3.71%'

$ws.Range("D39").Value = '''0.05627'
$ws.Range("E39").Value = '''-2.44%'

$ws.Range("D40").Value = '''0.007784'
$ws.Range("E40").Value = '''1.80%'

$ws.Range("D41").Value = '''0.1412'
$ws.Range("E41").Value = '''-1.12%'

$ws.Range("D42").Value = '''0.007324'

$ws.Range("D43").Value = '''0.002113'
$ws.Range("E43").Value = '''7.23%'

$ws.Range("D44").Value = '''0.007916'
$ws.Range("E44").Value = '''-1.42%'

$ws.Range("D45").Value = '''0.3509'

$ws.Range("D46").Value = '''0.00006957'
$ws.Range("E46").Value = '''-3.36%'

$ws.Range("E47").Value = '''0.16%'

$ws.Range("D48").Value = '''0.003490'
$ws.Range("E48").Value = '''0.06%'

$ws.Range("D49").Value = '''0.003532'
$ws.Range("E49").Value = '''39.94%'

$ws.Range("D50").Value = '''0.00002101'
$ws.Range("E50").Value = '''0.16%'

$ws.Range("D51").Value = '''0.0002001'
$ws.Range("E51").Value = '''0.16%'
